# Update the requirement list:
#  - Consolidate the three separate "Show name / Show image / Show description
#    of product ..." rows into a single "Show details of a product ..." row.
#  - Remove the two now-redundant rows.
#  - Renumber the remaining "Number" column so it stays sequential.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the three rows describing the catalogue click-through behaviour
# into a single row (row 6 currently holds "Show name of product...").
$ws.Range("C6").Value = "Show details of a product when clicking on it from the catalogue site"

# Row 7 (Working catalogue site...) stays where it is; rows 8 and 9
# ("Show image of product..." / "Show description of product...") are now
# redundant, so delete them - this shifts everything below up by two rows.
$ws.Rows("8:9").Delete()

# Renumber column B (the "Number" field) for the functional-requirements
# "A"/"B"/"C" rows so the sequence stays 1..24 after the two deleted rows.
$ws.Range("B8").Value = 4
$ws.Range("B9").Value = 5
$ws.Range("B10").Value = 6
$ws.Range("B11").Value = 7
$ws.Range("B12").Value = 8
$ws.Range("B13").Value = 9
$ws.Range("B14").Value = 10
$ws.Range("B15").Value = 11
$ws.Range("B16").Value = 12
$ws.Range("B17").Value = 13
$ws.Range("B18").Value = 14
$ws.Range("B19").Value = 15
$ws.Range("B20").Value = 16
$ws.Range("B21").Value = 17
$ws.Range("B22").Value = 18
$ws.Range("B23").Value = 19
$ws.Range("B24").Value = 20
$ws.Range("B25").Value = 21
$ws.Range("B26").Value = 22
$ws.Range("B27").Value = 23
$ws.Range("B28").Value = 24
